# Updated symbol list on Sat Jan  7 21:56:35 UTC 2023 with GitHub Actions
# Applies refreshed Price (D) / Volume(1h) (E) text values to the crypto table.

function Set-TextCellValue {
    param($Worksheet, $Row, $Col, $NewValue)
    $cell = $Worksheet.Cells.Item($Row, $Col)
    $cell.NumberFormat = "@"
    $cell.Value = $NewValue
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextCellValue $ws 2 4 "261.74"
Set-TextCellValue $ws 2 5 "1.05%"
Set-TextCellValue $ws 3 5 "0.49%"
Set-TextCellValue $ws 4 4 "4.713"
Set-TextCellValue $ws 4 5 "0.31%"
Set-TextCellValue $ws 5 4 "0.06210"
Set-TextCellValue $ws 5 5 "2.98%"
Set-TextCellValue $ws 6 4 "6.725"
Set-TextCellValue $ws 6 5 "0.71%"
Set-TextCellValue $ws 7 4 "0.8503"
Set-TextCellValue $ws 7 5 "-1.11%"
Set-TextCellValue $ws 8 4 "0.9097"
Set-TextCellValue $ws 8 5 "-1.28%"
Set-TextCellValue $ws 9 4 "0.1407"
Set-TextCellValue $ws 9 5 "0.80%"
Set-TextCellValue $ws 10 4 "0.04663"
Set-TextCellValue $ws 10 5 "-10.15%"
Set-TextCellValue $ws 11 4 "0.07089"
Set-TextCellValue $ws 11 5 "0.21%"
Set-TextCellValue $ws 12 4 "0.03165"
Set-TextCellValue $ws 12 5 "2.61%"
Set-TextCellValue $ws 13 4 "0.09056"
Set-TextCellValue $ws 13 5 "-0.90%"
Set-TextCellValue $ws 14 4 "0.001536"
Set-TextCellValue $ws 14 5 "0.04%"
Set-TextCellValue $ws 15 4 "0.0006151"
Set-TextCellValue $ws 15 5 "1.81%"
Set-TextCellValue $ws 16 4 "0.006128"
Set-TextCellValue $ws 16 5 "0.78%"
Set-TextCellValue $ws 17 4 "3.467"
Set-TextCellValue $ws 18 4 "3.168"
Set-TextCellValue $ws 18 5 "-0.17%"
Set-TextCellValue $ws 19 4 "2.177"
Set-TextCellValue $ws 21 4 "0.1299"
Set-TextCellValue $ws 21 5 "0.13%"
Set-TextCellValue $ws 22 4 "4.108"
Set-TextCellValue $ws 22 5 "-0.29%"
Set-TextCellValue $ws 23 5 "-0.49%"
Set-TextCellValue $ws 24 4 "0.001215"
Set-TextCellValue $ws 24 5 "-0.06%"
Set-TextCellValue $ws 25 4 "0.004137"
Set-TextCellValue $ws 25 5 "2.57%"
Set-TextCellValue $ws 27 4 "0.0001616"
Set-TextCellValue $ws 27 5 "6.09%"
Set-TextCellValue $ws 40 4 "0.03894"
Set-TextCellValue $ws 40 5 "1.19%"
Set-TextCellValue $ws 41 4 "0.1115"
Set-TextCellValue $ws 41 5 "0.03%"
Set-TextCellValue $ws 42 4 "0.004131"
Set-TextCellValue $ws 42 5 "2.62%"
Set-TextCellValue $ws 43 5 "-0.75%"
Set-TextCellValue $ws 44 5 "-11.63%"
Set-TextCellValue $ws 45 4 "0.00005172"
Set-TextCellValue $ws 45 5 "1.20%"
Set-TextCellValue $ws 46 5 "0.06%"
Set-TextCellValue $ws 47 4 "0.03590"
Set-TextCellValue $ws 47 5 "-34.19%"
Set-TextCellValue $ws 48 5 "23.17%"
Set-TextCellValue $ws 49 5 "0.06%"
Set-TextCellValue $ws 50 5 "0.06%"
